$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data ordering (Core Rules block moved to top, book/Roomlands block shifted down,
# Kyodaina/Dungeons block shifted down) per the restructured lookup table.
$data = @(
    @("Explore","Core Rules"),
    @("Ocean Explore","Core Rules"),
    @("Items","Core Rules"),
    @("Events","Core Rules"),
    @("Opponent Intentions","Core Rules"),
    @("Weapon Types","Core Rules"),
    @("Magic Types","Core Rules"),
    @("Body Types","Core Rules"),
    @("Rewards","Core Rules"),
    @("Combat Modifier","Core Rules"),
    @("Attack Types","Core Rules"),
    @("Default Oracle","Core Rules"),
    @("Two Player Oracle","Core Rules"),
    @("Complications","Core Rules"),
    @("Storytelling Incites","Core Rules"),
    @("Storytelling Subjects","Core Rules"),
    @("Storytelling Twists","Core Rules"),
    @("NPC Names","Core Rules"),
    @("NPC Looks","Core Rules"),
    @("NPC Characteristics","Core Rules"),
    @("Weather","Core Rules"),
    @("City Building","Core Rules"),
    @("Hunters Guild Quest Locations","Core Rules"),
    @("Hunters Guild Quest Distances","Core Rules"),
    @("Hunters Guild Quest Twists","Core Rules"),
    @("Hunters Guild Quest Rewards","Core Rules"),
    @("Rooklings","Core Rules"),
    @("Tundra Explore","Roomlands(Extended core)"),
    @("Tundra Vilalge","Roomlands(Extended core)"),
    @("Temple Skills","Roomlands(Extended core)"),
    @("Tundroom Characters","Roomlands(Extended core)"),
    @("Rooknaught Rewards","Roomlands(Extended core)"),
    @("Rafters encounters","Roomlands(Extended core)"),
    @("Crackways encounters","Roomlands(Extended core)"),
    @("Battlements Encounters","Roomlands(Extended core)"),
    @("Exposure event prompts","Roomlands(Extended core)"),
    @("Battlements attack type","Roomlands(Extended core)"),
    @("Biomes","Roomlands(Extended core)"),
    @("Rooks Rest Encounter","Roomlands(Extended core)"),
    @("Parapette Hunters Guild Quests","Roomlands(Extended core)"),
    @("Parapette Hunters Guild Quest Locations","Roomlands(Extended core)"),
    @("Parapette Hunters Guild Quest Distances","Roomlands(Extended core)"),
    @("Parapette Hunters Guild Quest Rewards","Roomlands(Extended core)"),
    @("The Market District","Roomlands(Extended core)"),
    @("The Crescents","Roomlands(Extended core)"),
    @("Guildtown","Roomlands(Extended core)"),
    @("The Lanes","Roomlands(Extended core)"),
    @("Rookhome Maintenance","Roomlands(Extended core)"),
    @("The Halls of Torrents","Roomlands(Extended core)"),
    @("The Halls of Winds","Roomlands(Extended core)"),
    @("The Halls of Gears","Roomlands(Extended core)"),
    @("Temple Skills (Stone Fist)","Kyodaina"),
    @("Kyodaina Zones","Kyodaina"),
    @("Kyodaina Explore","Kyodaina"),
    @("NPCs of the spirit forest","Kyodaina"),
    @("Village of Eda- Encounter table","Kyodaina"),
    @("Village of Eda- Event results table","Kyodaina"),
    @("Temple skills (Stone Face)","Kyodaina"),
    @("The Fangs of Fire NPCs","Kyodaina"),
    @("Shrine of Memory Key Locations","Kyodaina"),
    @("The Temple of the Stone Acolytes","Kyodaina"),
    @("The Hori Archipelago NPCs","Kyodaina"),
    @("Aanya's Prompt Table","Kyodaina"),
    @("The Temple of the Stone Dial","Kyodaina"),
    @("The Skylands NPCs","Kyodaina"),
    @("Sklyands Key Locations","Kyodaina"),
    @("The Temple of the Gargoyle","Kyodaina"),
    @("Shiro Standing prompts","Kyodaina"),
    @("Labyrinth Chambers","Dungeons"),
    @("Labyrinth Looks","Dungeons"),
    @("Labyrinth Items","Dungeons"),
    @("Labyrinth Guardians","Dungeons"),
    @("Labyrinth Guardians Rewards","Dungeons"),
    @("Dungeons Exploration","Dungeons"),
    @("Labyrinth Treasures","Dungeons"),
    @("The Oubliette NPCs","Dungeons"),
    @("Oubliette Exploration","Dungeons")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $ws.Cells.Item($i + 2, 1).Value = $row[0]
    $ws.Cells.Item($i + 2, 2).Value = $row[1]
}

# Row 39 (previously "Biomes" with a manually-applied fill format) no longer carries
# that leftover cell formatting after the reshuffle.
$ws.Range("A39").ClearFormats()

# Restore default selection to A2.
[void]$ws.Range("A2").Select()
